# "new object to align starting population with targets"
#
# The UK_cohabitation_adjustment and UK_fertility_adjustment sheets hold a
# year-by-year correction factor used to align the simulated starting
# population with observed targets. This replaces the old calibrated
# factors (years 2020-2069, rows 12-61) with 0 so the new alignment object
# starts from a neutral baseline.

$wb = $excel.ActiveWorkbook

$cohab = $wb.Worksheets.Item("UK_cohabitation_adjustment")
$fert  = $wb.Worksheets.Item("UK_fertility_adjustment")

# Zero-out the adjustment factors for years 2020-2069 (rows 12-61, column B)
$cohab.Range("B12:B61").Value = 0
$fert.Range("B12:B61").Value = 0

# Reflect the author's on-screen state when the file was saved: scrolled so
# row 24 is at the top, with B2:B61 selected on both sheets.
$cohab.Activate()
[void]$cohab.Range("B2:B61").Select()
$cohabWin = $excel.ActiveWindow
$cohabWin.ScrollRow = 24
$cohabWin.ScrollColumn = 1

$fert.Activate()
[void]$fert.Range("B2:B61").Select()
$fertWin = $excel.ActiveWindow
$fertWin.ScrollRow = 24
$fertWin.ScrollColumn = 1

# Restore the originally active sheet/tab (UK_fertility_adjustment)
$fert.Activate()
